# Rename the inline logo pictures embedded in the headers/footers.
#
# The document carries two distinct logos, each placed inline as the last
# run of the last paragraph of the header/footer part that hosts it:
#   - the Pearson logo (alt text contains "PearsonLogo") currently named
#     "image1.png" -> should become "image2.png"
#   - the BTec logo   (alt text "BTec_Logo-Orange")       currently named
#     "image2.jpg" -> should become "image1.jpg"
#
# This swap happens in all three places the logos appear: the default
# footer, the first-page footer and the first-page header.

$d = $word.ActiveDocument

function Rename-LogoInHeaderFooter($hf) {
    if (-not $hf.Exists) {
        return
    }
    $rng = $hf.Range
    $paraCount = $rng.Paragraphs.Count
    if ($paraCount -lt 1) {
        return
    }
    # The picture run lives in the last paragraph of the part; addressing
    # the inline shape through that specific paragraph (rather than the
    # whole-story Range) is what keeps the write from bouncing off a
    # stale-handle error when the part has several paragraphs.
    $p = $rng.Paragraphs($paraCount)
    if ($p.Range.InlineShapes.Count -lt 1) {
        return
    }

    $shp = $p.Range.InlineShapes(1)
    $alt = $shp.AlternativeText

    if ($alt -like "*PearsonLogo*") {
        $shp.Name = "image2.png"
    } elseif ($alt -like "*BTec_Logo-Orange*") {
        $shp.Name = "image1.jpg"
    }
}

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        Rename-LogoInHeaderFooter $sec.Headers($i)
        Rename-LogoInHeaderFooter $sec.Footers($i)
    }
}

Write-Host "Logo inline shapes renamed."
